$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Numeric corrections in rows 9 and 10 (columns F:L)
# ---------------------------------------------------------------------------
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 20
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 30
$ws.Range("L9").Value = 10

$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 20
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 10
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 30
$ws.Range("L10").Value = 10

# ---------------------------------------------------------------------------
# 2) Change the text "0.1" -> "0.4" in a handful of scattered cells.
#
# A plain `.Value = "0.4"` would be auto-coerced to the NUMBER 0.4 by
# Excel's normal "smart" entry parsing (exactly like typing 0.4 into a
# cell on a real sheet), which is not what we want here: the source file
# stores these as literal TEXT (shared string), not as a number, and the
# cell's existing style/format must stay exactly as-is.
#
# The reliable way to push literal text that merely *looks* like a number
# into a cell without Excel re-interpreting it, and without touching the
# destination cell's formatting, is to stage the text in a scratch cell
# that is explicitly formatted as Text, copy it, and PasteSpecial only the
# VALUES into each destination (xlPasteValues leaves the destination's own
# number format/style untouched).
# ---------------------------------------------------------------------------
$scratch = $ws.Range("AB100")
$scratch.NumberFormat = "@"
$scratch.Value = "0.4"
$scratch.Copy()

$targets = @("S8", "V9", "X10", "Z11", "S12", "S13", "S14", "S15", "S16", "S17", "S18", "S19", "S20")
foreach ($t in $targets) {
    $ws.Range($t).PasteSpecial(-4163)
}

# Clean up the scratch cell so it doesn't leave any trace behind.
$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Restore the current selection shown in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("P9:R9").Select()
